# Generate Report for Handback
#
# The handback transform for e2e\1f586067-53a6-4d94-b147-c3cd3ea52e96.md failed
# (file-name mismatch between the handback package and the original handoff)
# for both the zh-cn and de-de targets. Update the status + error-detail
# columns on each locale sheet (and the rolled-up status on the Overview
# sheet) to reflect the failure, and widen the "Error Detail" column so the
# new message is readable.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Overview sheet: per-locale status columns (zh-cn, de-de) for the
# 1f586067-... row (row 3).
$ws1.Range("E3").Value = $newStatus
$ws1.Range("F3").Value = $newStatus

# Locale sheets: "Status" column (C) for the same row.
$ws2.Range("C3").Value = $newStatus
$ws3.Range("C3").Value = $newStatus

# Locale sheets: "Error Detail" column (P) for the same row, explaining why
# the handback transform failed.
$ws2.Range("P3").Value = "Handback file name: p5gru1ny.xtf is different with handoff file name: 1f586067-53a6-4d94-b147-c3cd3ea52e96.a883ed08c88991394c7ea2589578b01d86288442.zh-cn."
$ws3.Range("P3").Value = "Handback file name: p5gru1ny.xtf is different with handoff file name: 1f586067-53a6-4d94-b147-c3cd3ea52e96.a883ed08c88991394c7ea2589578b01d86288442.de-de."

# Widen the "Error Detail" column (16th / P) on both locale sheets so the
# new, longer message is visible (stored column width 40 <=> ColumnWidth 39.17).
$ws2.Columns.Item(16).ColumnWidth = 39.17
$ws3.Columns.Item(16).ColumnWidth = 39.17
